$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New content rows (12-15) ---

# Row 12: plain heading text
$ws.Range("B12").Value = "Adatgyűjtési hatékonyság összehasonlítása (Scrapy vs BeautifulSoup)"

# Row 13: rich text "Scrapy és a BeautifulSoup az adatgyűjtés technológiája."
# with "Scrapy és a BeautifulSoup" bold (inherited cell font), "adatgyűjtés" italic
$b13text = "Scrapy és a BeautifulSoup az adatgyűjtés technológiája."
$ws.Range("B13").Value = $b13text
$ws.Range("B13").Font.Bold = $true

$italicStart = $b13text.IndexOf("adatgyűjtés") + 1
$italicLen = "adatgyűjtés".Length
$afterItalicStart = $italicStart + $italicLen
$afterItalicLen = $b13text.Length - ($afterItalicStart - 1)

# " az " (between "BeautifulSoup" and "adatgyűjtés") -> regular weight
$midStart = "Scrapy és a BeautifulSoup".Length + 1
$midLen = $italicStart - $midStart
$ws.Range("B13").Characters($midStart, $midLen).Font.Bold = $false

# "adatgyűjtés" -> italic, not bold
$chars = $ws.Range("B13").Characters($italicStart, $italicLen)
$chars.Font.Italic = $true
$chars.Font.Bold = $false

# " technológiája." -> regular weight
$ws.Range("B13").Characters($afterItalicStart, $afterItalicLen).Font.Bold = $false

# Row 14: plain paragraph text
$ws.Range("B14").Value = "Az adatgyűjtési fázis során a Scrapy és a BeautifulSoup könyvtárak teljesítményét hasonlítottam össze különböző URL-ek esetén, mérve a feldolgozási időt,adatgyűjtései sebességet, memóriahasználatot"

# Row 15: rich text with a bold middle portion
$boldPhrase = "a HTML elemek manuális feldolgozása jelentős fejlesztői munkát igényel"
$b15text = "A dolgozat célja a Scrapy és a BeautifulSoup könyvtárak összehasonlítása volt az URL-ekből történő adatgyűjtés hatékonysága szempontjából. A vizsgálat során kiderült, hogy bár mindkét eszköz képes strukturált adatok lekérésére, " + $boldPhrase + ", különösen eltérő oldalstruktúrák esetén."
$ws.Range("B15").Value = $b15text

$boldStart = $b15text.IndexOf($boldPhrase) + 1
$boldLen = $boldPhrase.Length
$ws.Range("B15").Characters($boldStart, $boldLen).Font.Bold = $true

$ws.Range("H10").Select()
